$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '26.413.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.95%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '1.671.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +1.08%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.58%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '221.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +1.61%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '0.5333'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +0.72%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.51%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.2662'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +1.57%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '0.06383'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +0.98%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '20.85'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +2.14%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.07856'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +0.73%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '4.531'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +0.31%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '1.680.35'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +0.71%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '1.902.33'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '0.5610'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +2.24%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '0.0₅8192'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +0.15%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '66.08'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +1.12%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '26.431.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +1.12%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +0.56%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '4.724'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +2.72%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '198.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +3.74%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '10.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +2.16%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '6.072'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +1.21%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '1.012'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.53%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '146.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +0.63%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '0.1225'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +0.15%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '7.256'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +0.72%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '16.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +1.39%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '1.506'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +2.39%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '0.05917'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +3.54%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '1.289'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +1.29%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '3.560'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +0.35%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '3.331'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +2.14%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '1.611'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +1.38%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '0.9686'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +2.18%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '2.838'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +1.24%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '2.440'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +0.77%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '0.5840'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +2.04%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '0.01616'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +0.58%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '1.077.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +3.71%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '5.934'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '0.8651'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +1.57%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +0.58%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '103.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -0.71%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '1.812.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +1.02%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '58.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +3.23%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +1.63%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +0.67%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '0.4420'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +1.49%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '8.035'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +2.64%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '0.05164'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.18%  '
$ws.Range('E51').Style = 'Normal'
